# Insert 3 new price rows for "June Pearl" ($/bandeja 18 kilos granel)
# right before the existing "Nectar Crest" block (old row 176), pushing
# every following row down by 3 — which also appends 3 rows at the very
# bottom of the used range (old last 3 rows -> new rows 282-284).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("176:178").Insert()

# Shared (constant) column values for this sheet's data rows
$mercadoId = 11
$mercado   = "Vega Monumental Concepción"
$region    = "Bíobío"
$codreg    = 8
$tipo      = "Fruta"
$prodId    = 100103
$producto  = "Frutos de hueso (carozo)"
$catId     = 100103006
$categoria = "Nectarín"
$origen    = "Región de O'Higgins"

# Row 176: June Pearl / Especial
$ws.Cells.Item(176,1).Value  = $mercadoId
$ws.Cells.Item(176,2).Value  = $mercado
$ws.Cells.Item(176,3).Value  = $region
$ws.Cells.Item(176,4).Value2 = 44603
$ws.Cells.Item(176,5).Value  = $codreg
$ws.Cells.Item(176,6).Value  = $tipo
$ws.Cells.Item(176,7).Value  = $prodId
$ws.Cells.Item(176,8).Value  = $producto
$ws.Cells.Item(176,9).Value  = $catId
$ws.Cells.Item(176,10).Value = $categoria
$ws.Cells.Item(176,11).Value = "June Pearl"
$ws.Cells.Item(176,12).Value = "Especial"
$ws.Cells.Item(176,13).Value = 50
$ws.Cells.Item(176,14).Value = 15000
$ws.Cells.Item(176,15).Value = 15000
$ws.Cells.Item(176,16).Value = 15000
$ws.Cells.Item(176,17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(176,18).Value = $origen
$ws.Cells.Item(176,19).Value = 833
$ws.Cells.Item(176,20).Value = 18

# Row 177: June Pearl / Primera
$ws.Cells.Item(177,1).Value  = $mercadoId
$ws.Cells.Item(177,2).Value  = $mercado
$ws.Cells.Item(177,3).Value  = $region
$ws.Cells.Item(177,4).Value2 = 44603
$ws.Cells.Item(177,5).Value  = $codreg
$ws.Cells.Item(177,6).Value  = $tipo
$ws.Cells.Item(177,7).Value  = $prodId
$ws.Cells.Item(177,8).Value  = $producto
$ws.Cells.Item(177,9).Value  = $catId
$ws.Cells.Item(177,10).Value = $categoria
$ws.Cells.Item(177,11).Value = "June Pearl"
$ws.Cells.Item(177,12).Value = "Primera"
$ws.Cells.Item(177,13).Value = 100
$ws.Cells.Item(177,14).Value = 13000
$ws.Cells.Item(177,15).Value = 13000
$ws.Cells.Item(177,16).Value = 13000
$ws.Cells.Item(177,17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(177,18).Value = $origen
$ws.Cells.Item(177,19).Value = 722
$ws.Cells.Item(177,20).Value = 18

# Row 178: June Pearl / Segunda
$ws.Cells.Item(178,1).Value  = $mercadoId
$ws.Cells.Item(178,2).Value  = $mercado
$ws.Cells.Item(178,3).Value  = $region
$ws.Cells.Item(178,4).Value2 = 44603
$ws.Cells.Item(178,5).Value  = $codreg
$ws.Cells.Item(178,6).Value  = $tipo
$ws.Cells.Item(178,7).Value  = $prodId
$ws.Cells.Item(178,8).Value  = $producto
$ws.Cells.Item(178,9).Value  = $catId
$ws.Cells.Item(178,10).Value = $categoria
$ws.Cells.Item(178,11).Value = "June Pearl"
$ws.Cells.Item(178,12).Value = "Segunda"
$ws.Cells.Item(178,13).Value = 100
$ws.Cells.Item(178,14).Value = 11000
$ws.Cells.Item(178,15).Value = 11000
$ws.Cells.Item(178,16).Value = 11000
$ws.Cells.Item(178,17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(178,18).Value = $origen
$ws.Cells.Item(178,19).Value = 611
$ws.Cells.Item(178,20).Value = 18
